$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.030.48'
$ws.Range('E2').Value = '  -4.55%  '
$ws.Range('D3').Value = '2.500.88'
$ws.Range('E3').Value = '  -2.94%  '
$styleSave = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $styleSave
$ws.Range('E4').Value = '  -0.18%  '
$styleSave = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.80'
$ws.Range('D5').Style = $styleSave
$ws.Range('E5').Value = '  -2.27%  '
$styleSave = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.09'
$ws.Range('D6').Style = $styleSave
$ws.Range('E6').Value = '  -6.56%  '
$styleSave = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.995'
$ws.Range('D7').Style = $styleSave
$ws.Range('E7').Value = '  -0.48%  '
$styleSave = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('D8').Style = $styleSave
$ws.Range('E8').Value = '  -2.91%  '
$ws.Range('D9').Value = '2.533.10'
$ws.Range('E9').Value = '  -1.90%  '
$styleSave = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0998'
$ws.Range('D10').Style = $styleSave
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('E11').Value = '  -2.54%  '
$styleSave = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.45'
$ws.Range('D12').Style = $styleSave
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('D14').Value = '2.936.73'
$ws.Range('E14').Value = '  -3.20%  '
$styleSave = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.76'
$ws.Range('D15').Style = $styleSave
$ws.Range('E15').Value = '  -6.28%  '
$ws.Range('D16').Value = '58.858.26'
$ws.Range('E16').Value = '  -4.77%  '
$ws.Range('E17').Value = '  -3.08%  '
$ws.Range('D18').Value = '2.511.36'
$ws.Range('E18').Value = '  -2.73%  '
$styleSave = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.31'
$ws.Range('D19').Style = $styleSave
$ws.Range('E19').Value = '  -2.16%  '
$styleSave = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.28'
$ws.Range('D20').Style = $styleSave
$ws.Range('E20').Value = '  -5.22%  '
$styleSave = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.02'
$ws.Range('D21').Style = $styleSave
$ws.Range('E21').Value = '  -4.34%  '
$styleSave = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = $styleSave
$ws.Range('E22').Value = '  +0.16%  '
$styleSave = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.76'
$ws.Range('D23').Style = $styleSave
$ws.Range('E23').Value = '  -4.04%  '
$styleSave = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.82'
$ws.Range('D24').Style = $styleSave
$ws.Range('E24').Value = '  -2.60%  '
$styleSave = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.439'
$ws.Range('D25').Style = $styleSave
$ws.Range('E25').Value = '  -10.30%  '
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('D27').Value = '2.612.08'
$ws.Range('E27').Value = '  -3.15%  '
$styleSave = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.994'
$ws.Range('D28').Style = $styleSave
$ws.Range('E28').Value = '  -0.55%  '
$styleSave = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.76'
$ws.Range('D29').Style = $styleSave
$ws.Range('E29').Value = '  -4.31%  '
$styleSave = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.78'
$ws.Range('D30').Style = $styleSave
$ws.Range('E30').Value = '  -5.35%  '
$ws.Range('D31').Value = '0.0₃0779'
$ws.Range('E31').Value = '  -6.27%  '
$ws.Range('E32').Value = '  -5.00%  '
$styleSave = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('D33').Style = $styleSave
$ws.Range('E33').Value = '  -8.54%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$styleSave = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.996'
$ws.Range('D34').Style = $styleSave
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$styleSave = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '159.33'
$ws.Range('D35').Style = $styleSave
$ws.Range('E35').Value = '  -1.99%  '
$ws.Range('E36').Value = '  +4.58%  '
$styleSave = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.58'
$ws.Range('D37').Style = $styleSave
$ws.Range('E37').Value = '  -2.67%  '
$styleSave = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.39'
$ws.Range('D38').Style = $styleSave
$ws.Range('E38').Value = '  -9.07%  '
$ws.Range('E39').Value = '  -8.60%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$styleSave = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.71'
$ws.Range('D40').Style = $styleSave
$ws.Range('E40').Value = '  -4.85%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$styleSave = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '306.19'
$ws.Range('D41').Style = $styleSave
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$styleSave = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.83'
$ws.Range('D42').Style = $styleSave
$ws.Range('E42').Value = '  -1.72%  '
$styleSave = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.823'
$ws.Range('D43').Style = $styleSave
$ws.Range('E43').Value = '  -8.12%  '
$ws.Range('E44').Value = '  -6.26%  '
$styleSave = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.993'
$ws.Range('D45').Style = $styleSave
$ws.Range('E45').Value = '  -0.55%  '
$styleSave = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.599'
$ws.Range('D46').Style = $styleSave
$ws.Range('E46').Value = '  -0.95%  '
$styleSave = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.77'
$ws.Range('D47').Style = $styleSave
$ws.Range('E47').Value = '  -1.47%  '
$styleSave = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.35'
$ws.Range('D48').Style = $styleSave
$ws.Range('E48').Value = '  +2.52%  '
$styleSave = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0930'
$ws.Range('D49').Style = $styleSave
$ws.Range('E49').Value = '  -3.45%  '
$styleSave = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.60'
$ws.Range('D50').Style = $styleSave
$ws.Range('E50').Value = '  -4.59%  '
$styleSave = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0518'
$ws.Range('D51').Style = $styleSave
$ws.Range('E51').Value = '  -4.72%  '
